$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings; force Text format so
# Excel does not silently coerce them to floating-point numbers (which
# would lose formatting like trailing zeros and introduce fp noise).
$dCells = @("D2","D3","D5","D6","D8","D9","D11","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D27","D28","D29","D30","D32","D34","D36","D39","D40","D41","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.208.10"
$ws.Range("D3").Value = "2.343.30"
$ws.Range("D5").Value = "545.67"
$ws.Range("D6").Value = "132.27"
$ws.Range("D8").Value = "0.586"
$ws.Range("D9").Value = "2.341.69"
$ws.Range("D11").Value = "5.52"
$ws.Range("D14").Value = "23.86"
$ws.Range("D15").Value = "2.760.29"
$ws.Range("D16").Value = "60.234.93"
$ws.Range("D18").Value = "2.354.04"
$ws.Range("D19").Value = "10.61"
$ws.Range("D20").Value = "4.15"
$ws.Range("D22").Value = "314.44"
$ws.Range("D23").Value = "0.998"
$ws.Range("D24").Value = "63.54"
$ws.Range("D25").Value = "0.172"
$ws.Range("D27").Value = "7.91"
$ws.Range("D28").Value = "1.36"
$ws.Range("D29").Value = "1.76"
$ws.Range("D30").Value = "171.65"
$ws.Range("D32").Value = "0.0₃0729"
$ws.Range("D34").Value = "1.40"
$ws.Range("D36").Value = "18.06"
$ws.Range("D39").Value = "4.16"
$ws.Range("D40").Value = "324.25"
$ws.Range("D41").Value = "38.10"
$ws.Range("D43").Value = "141.28"
$ws.Range("D44").Value = "3.46"
$ws.Range("D45").Value = "0.0947"
$ws.Range("D46").Value = "19.56"
$ws.Range("D47").Value = "0.0498"
$ws.Range("D48").Value = "0.560"
$ws.Range("D50").Value = "11.02"
$ws.Range("D51").Value = "0.0₆0206"

# Restore default (unstyled) look now that the text is stored safely.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Column E values are already non-numeric (leading/trailing spaces + % sign)
# so a plain Value assignment keeps them as text.
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  +6.99%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +8.10%  "
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  +13.43%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("E34").Value = "  +14.53%  "
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("E40").Value = "  +13.50%  "
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +9.47%  "
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +13.03%  "
